{"js": "// Update the answer table: each cell's text is a \"NNN\u00f7N=NN, N\" style\n// division fact. Replace the old fact text with the new one, cell by cell,\n// using a plain text search + replace so formatting (fonts/size) on the\n// run is preserved.\nconst replacements = [\n  [\"781\u00f73=260, 1\", \"721\u00f79=80, 1\"],\n  [\"737\u00f77=105, 2\", \"258\u00f74=64, 2\"],\n  [\"663\u00f79=73, 6\", \"996\u00f72=498, 0\"],\n  [\"441\u00f79=49, 0\", \"110\u00f73=36, 2\"],\n  [\"469\u00f77=67, 0\", \"149\u00f73=49, 2\"],\n  [\"708\u00f77=101, 1\", \"974\u00f78=121, 6\"],\n  [\"862\u00f77=123, 1\", \"387\u00f77=55, 2\"],\n  [\"173\u00f76=28, 5\", \"178\u00f75=35, 3\"],\n  [\"271\u00f78=33, 7\", \"268\u00f76=44, 4\"],\n  [\"567\u00f77=81, 0\", \"461\u00f72=230, 1\"],\n  [\"360\u00f73=120, 0\", \"692\u00f75=138, 2\"],\n  [\"835\u00f79=92, 7\", \"377\u00f72=188, 1\"],\n  [\"833\u00f79=92, 5\", \"101\u00f74=25, 1\"],\n  [\"749\u00f78=93, 5\", \"385\u00f78=48, 1\"],\n  [\"105\u00f79=11, 6\", \"759\u00f77=108, 3\"],\n  [\"507\u00f76=84, 3\", \"773\u00f77=110, 3\"],\n  [\"384\u00f72=192, 0\", \"101\u00f72=50, 1\"],\n  [\"554\u00f74=138, 2\", \"583\u00f79=64, 7\"],\n  [\"860\u00f74=215, 0\", \"558\u00f73=186, 0\"],\n  [\"484\u00f74=121, 0\", \"979\u00f74=244, 3\"],\n  [\"794\u00f76=132, 2\", \"467\u00f76=77, 5\"],\n  [\"716\u00f76=119, 2\", \"926\u00f75=185, 1\"],\n  [\"712\u00f73=237, 1\", \"554\u00f78=69, 2\"],\n  [\"769\u00f73=256, 1\", \"968\u00f74=242, 0\"],\n  [\"179\u00f78=22, 3\", \"405\u00f79=45, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the answer table: each cell's text is a \"NNN\u00f7N=NN, N\" style\n# division fact. Replace the old fact text with the new one, cell by cell,\n# using Find/Replace on the document's Content range (format-preserving:\n# only the run's text changes).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"781\u00f73=260, 1\", \"721\u00f79=80, 1\"),\n  @(\"737\u00f77=105, 2\", \"258\u00f74=64, 2\"),\n  @(\"663\u00f79=73, 6\", \"996\u00f72=498, 0\"),\n  @(\"441\u00f79=49, 0\", \"110\u00f73=36, 2\"),\n  @(\"469\u00f77=67, 0\", \"149\u00f73=49, 2\"),\n  @(\"708\u00f77=101, 1\", \"974\u00f78=121, 6\"),\n  @(\"862\u00f77=123, 1\", \"387\u00f77=55, 2\"),\n  @(\"173\u00f76=28, 5\", \"178\u00f75=35, 3\"),\n  @(\"271\u00f78=33, 7\", \"268\u00f76=44, 4\"),\n  @(\"567\u00f77=81, 0\", \"461\u00f72=230, 1\"),\n  @(\"360\u00f73=120, 0\", \"692\u00f75=138, 2\"),\n  @(\"835\u00f79=92, 7\", \"377\u00f72=188, 1\"),\n  @(\"833\u00f79=92, 5\", \"101\u00f74=25, 1\"),\n  @(\"749\u00f78=93, 5\", \"385\u00f78=48, 1\"),\n  @(\"105\u00f79=11, 6\", \"759\u00f77=108, 3\"),\n  @(\"507\u00f76=84, 3\", \"773\u00f77=110, 3\"),\n  @(\"384\u00f72=192, 0\", \"101\u00f72=50, 1\"),\n  @(\"554\u00f74=138, 2\", \"583\u00f79=64, 7\"),\n  @(\"860\u00f74=215, 0\", \"558\u00f73=186, 0\"),\n  @(\"484\u00f74=121, 0\", \"979\u00f74=244, 3\"),\n  @(\"794\u00f76=132, 2\", \"467\u00f76=77, 5\"),\n  @(\"716\u00f76=119, 2\", \"926\u00f75=185, 1\"),\n  @(\"712\u00f73=237, 1\", \"554\u00f78=69, 2\"),\n  @(\"769\u00f73=256, 1\", \"968\u00f74=242, 0\"),\n  @(\"179\u00f78=22, 3\", \"405\u00f79=45, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
